# "Generate Report for Handback"
#
# The workbook tracks localization handoff/handback status for two
# languages (zh-cn, de-de). This change reports the handback: it marks
# the Overview status text as handed-back, fills in the "Latest Target
# File" / "Latest Handback File" columns (F/G) with the relevant
# md / xlf file links (mirroring the existing source-file / handoff-file
# hyperlinks already present in columns A/D), and stamps the "Latest
# Handback DateTime" column (H) with the real handback timestamp where
# applicable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet: status text "Ready for handoff" -> "Handed back: in
#    sync with en-US" for both language columns / both file rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$handedBack = "Handed back: in sync with en-US"
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------
# Helper: add a hyperlink cell (value + display text + style) without
# clobbering the rest of the row.
# ---------------------------------------------------------------------
function Add-ReportLink($ws, $cellRef, $address, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, [System.Type]::Missing, [System.Type]::Missing, $text) | Out-Null
}

# ---------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Row 2 - 71ecb103-57f2-45c9-9a02-c6f5cbed3946
Add-ReportLink $zhcn "F2" "https://github.com/OpenLocalizationTest/oltest/blob/a28127b7c67e1dcd8632f2663e4c55cabf70a205/e2e/71ecb103-57f2-45c9-9a02-c6f5cbed3946.md" "71ecb103-57f2-45c9-9a02-c6f5cbed3946.md"
Add-ReportLink $zhcn "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e7d20289af759a612f3cf7b68c0ef3b873a6900/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/71ecb103-57f2-45c9-9a02-c6f5cbed3946.bb27e8fa41549806020f54da1da8248823d0c53e.zh-cn.xlf" "71ecb103-57f2-45c9-9a02-c6f5cbed3946.bb27e8fa41549806020f54da1da8248823d0c53e.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-03-11 20:10:38"

# Row 3 - 904dc162-2b0e-457a-aa64-2894a933e14a (no handback datetime yet)
Add-ReportLink $zhcn "F3" "https://github.com/OpenLocalizationTest/oltest/blob/a28127b7c67e1dcd8632f2663e4c55cabf70a205/e2e/904dc162-2b0e-457a-aa64-2894a933e14a.md" "904dc162-2b0e-457a-aa64-2894a933e14a.md"
Add-ReportLink $zhcn "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e7d20289af759a612f3cf7b68c0ef3b873a6900/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/904dc162-2b0e-457a-aa64-2894a933e14a.e7982190315cd4f6d7cea62a4363d8ec68b45715.zh-cn.xlf" "904dc162-2b0e-457a-aa64-2894a933e14a.e7982190315cd4f6d7cea62a4363d8ec68b45715.zh-cn.xlf"

# ---------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Row 2 - 71ecb103-57f2-45c9-9a02-c6f5cbed3946
Add-ReportLink $dede "F2" "https://github.com/OpenLocalizationTest/oltest/blob/a28127b7c67e1dcd8632f2663e4c55cabf70a205/e2e/71ecb103-57f2-45c9-9a02-c6f5cbed3946.md" "71ecb103-57f2-45c9-9a02-c6f5cbed3946.md"
Add-ReportLink $dede "G2" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18bd34f6c2225e617c5c173e67209e6166ffe964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/71ecb103-57f2-45c9-9a02-c6f5cbed3946.bb27e8fa41549806020f54da1da8248823d0c53e.de-de.xlf" "71ecb103-57f2-45c9-9a02-c6f5cbed3946.bb27e8fa41549806020f54da1da8248823d0c53e.de-de.xlf"
$dede.Range("H2").Value = "2016-03-11 20:10:46"

# Row 3 - 904dc162-2b0e-457a-aa64-2894a933e14a
Add-ReportLink $dede "F3" "https://github.com/OpenLocalizationTest/oltest/blob/a28127b7c67e1dcd8632f2663e4c55cabf70a205/e2e/904dc162-2b0e-457a-aa64-2894a933e14a.md" "904dc162-2b0e-457a-aa64-2894a933e14a.md"
Add-ReportLink $dede "G3" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18bd34f6c2225e617c5c173e67209e6166ffe964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/904dc162-2b0e-457a-aa64-2894a933e14a.e7982190315cd4f6d7cea62a4363d8ec68b45715.de-de.xlf" "904dc162-2b0e-457a-aa64-2894a933e14a.e7982190315cd4f6d7cea62a4363d8ec68b45715.de-de.xlf"
$dede.Range("H3").Value = "2016-03-11 20:10:46"

Write-Host "Handback report generated."
